# "filling marker info into status 8"
# Fill the marker_1 column (J) with "NAT" for the rows that were still
# missing it (everything except the untagged parent-strain rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(4, 5, 6, 7, 8, 9, 10, 22, 23, 26, 27, 28, 29, 30, 31, 32, 44, 45)

foreach ($r in $rows) {
    $ws.Range("J$r").Value = "NAT"
}

# Update the active selection to match the saved view state.
$ws.Range("Q20").Select()
